$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in the checklist (column B) marked "Y" (Yes)
$yRows = @(38, 39, 40, 44, 45, 46, 47, 48, 49, 50, 51, 52, 60, 61, 62, 63, 64, 68, 69, 70, 71, 72, 77, 78, 80, 84, 88, 89, 95, 99, 101, 107, 108, 109, 110, 111, 112, 117, 118, 120, 122, 151, 152, 156, 157, 158, 160, 161, 162, 164, 174, 175, 176, 177, 181, 182, 183, 184, 185, 186, 187, 188, 189)

# Rows in the checklist (column B) marked "N" (No)
$nRows = @(41, 53, 54, 55, 190)

foreach ($r in $yRows) {
    $ws.Cells.Item($r, 2).Value = "Y"
}

foreach ($r in $nRows) {
    $ws.Cells.Item($r, 2).Value = "N"
}

# B74 keeps no value, but loses its red-box "fillable" formatting in favour of
# a plain, unlocked cell (the marker could not tick this particular item).
$ws.Range("B74").Style = "Normal"
$ws.Range("B74").Locked = $false

# The sheet view was scrolled down and the active selection moved while the
# marker was working through the list.
$excel.ActiveWindow.ScrollRow = 56
$ws.Range("B73").Select() | Out-Null
